$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr><w:r><w:t>This is an another line of Document where the actual line was written from Labotop' + [char]0x2026 + '</w:t></w:r><w:r><w:br/><w:t>By Karunakar</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$endRange.InsertXML($xml)
